$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The company in row 2 (Grand Industrial) has a new / updated website.
# Replace the old grandindustrial.com URL with the Oerlikon Fairfield page.
$ws.Range("C2").Value = "https://www.oerlikon.com/fairfield/en/"

# Move the active selection to C7 (where the cursor ended up after the edit).
[void]$ws.Range("C7").Select()

$wb.Save()
